# Updated symbol list on Tue Jan 17 17:50:30 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns on the crypto
# ticker sheet with the latest scraped snapshot. Values are stored as
# plain text (as the source scraper always wrote them), so each cell's
# NumberFormat is forced to Text ("@") before the write to stop Excel
# from silently reinterpreting look-alike numbers/percentages as native
# numeric values; the cell style is then reset to "Normal" so no stray
# formatting is left behind beyond the text content itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "301.96" }
    @{ Cell = "E2";  Value = "1.26%" }
    @{ Cell = "D3";  Value = "31.83" }
    @{ Cell = "E3";  Value = "0.48%" }
    @{ Cell = "D4";  Value = "5.022" }
    @{ Cell = "E4";  Value = "-0.94%" }
    @{ Cell = "E5";  Value = "-3.72%" }
    @{ Cell = "D6";  Value = "2.038" }
    @{ Cell = "E6";  Value = "-22.50%" }
    @{ Cell = "D7";  Value = "7.813" }
    @{ Cell = "E7";  Value = "0.01%" }
    @{ Cell = "D8";  Value = "3.776" }
    @{ Cell = "E8";  Value = "-1.20%" }
    @{ Cell = "D9";  Value = "0.9189" }
    @{ Cell = "E9";  Value = "-0.89%" }
    @{ Cell = "D10"; Value = "0.1743" }
    @{ Cell = "E10"; Value = "-1.06%" }
    @{ Cell = "D11"; Value = "0.07875" }
    @{ Cell = "E11"; Value = "5.43%" }
    @{ Cell = "D12"; Value = "0.08762" }
    @{ Cell = "E12"; Value = "-0.92%" }
    @{ Cell = "D13"; Value = "0.03099" }
    @{ Cell = "E13"; Value = "2.21%" }
    @{ Cell = "E14"; Value = "-0.07%" }
    @{ Cell = "D15"; Value = "0.001505" }
    @{ Cell = "E15"; Value = "-0.36%" }
    @{ Cell = "D16"; Value = "0.005829" }
    @{ Cell = "E16"; Value = "-0.90%" }
    @{ Cell = "D17"; Value = "3.463" }
    @{ Cell = "E17"; Value = "-2.79%" }
    @{ Cell = "D18"; Value = "2.268" }
    @{ Cell = "D19"; Value = "0.3270" }
    @{ Cell = "E19"; Value = "1.26%" }
    @{ Cell = "D20"; Value = "0.1290" }
    @{ Cell = "E20"; Value = "-2.10%" }
    @{ Cell = "D21"; Value = "4.156" }
    @{ Cell = "E21"; Value = "4.86%" }
    @{ Cell = "D23"; Value = "0.04593" }
    @{ Cell = "E23"; Value = "-0.17%" }
    @{ Cell = "D24"; Value = "0.001240" }
    @{ Cell = "E24"; Value = "-0.20%" }
    @{ Cell = "D25"; Value = "0.004475" }
    @{ Cell = "E25"; Value = "0.43%" }
    @{ Cell = "E26"; Value = "4.48%" }
    @{ Cell = "D39"; Value = "0.01735" }
    @{ Cell = "E39"; Value = "-1.83%" }
    @{ Cell = "D40"; Value = "0.04763" }
    @{ Cell = "E40"; Value = "5.15%" }
    @{ Cell = "D41"; Value = "0.007296" }
    @{ Cell = "E41"; Value = "6.65%" }
    @{ Cell = "E42"; Value = "0.40%" }
    @{ Cell = "D43"; Value = "0.002140" }
    @{ Cell = "E43"; Value = "0.16%" }
    @{ Cell = "D44"; Value = "0.01082" }
    @{ Cell = "E44"; Value = "9.81%" }
    @{ Cell = "D45"; Value = "0.00006068" }
    @{ Cell = "E45"; Value = "-6.25%" }
    @{ Cell = "E46"; Value = "0.15%" }
    @{ Cell = "D47"; Value = "0.003500" }
    @{ Cell = "E47"; Value = "-59.94%" }
    @{ Cell = "D48"; Value = "0.8204" }
    @{ Cell = "E48"; Value = "-0.02%" }
    @{ Cell = "E49"; Value = "0.15%" }
    @{ Cell = "E50"; Value = "0.15%" }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
    $cell.Style = "Normal"
}
